$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.217414975166321
$ws.Range("B1").Value = 2.228760957717896
$ws.Range("D1").Value = 1.433882236480713
$ws.Range("E1").Value = 0.8967984914779663
